# "Data Copy understands byte arrays. Data collected"
#
# 1. Add a new "DATA COPY RESULTS" sheet (a timing-results table shaped just
#    like the existing "MEM MAP RESULTS" sheet) at the end of the workbook,
#    fill it with the newly collected byte-array timings, and restyle the
#    average row to match the other result sheets.
# 2. Re-apply the bold/centered "average row" style to the B:E cells of the
#    two sheets whose average-row formulas were touched by the same pass
#    (JAVA SOCKETS RESULTS, MEM MAP RESULTS).
# 3. Leave the selection/active-tab state the way the author left it:
#    "1 Message Sent" ends up the active sheet (cell K15 selected), while
#    JAVA SOCKETS RESULTS and MEM MAP RESULTS keep their own in-sheet
#    selections but are no longer the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "DATA COPY RESULTS" sheet.
# ---------------------------------------------------------------------
$memMap = $wb.Worksheets.Item("MEM MAP RESULTS")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$memMap.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "DATA COPY RESULTS"

$ws.Range("A1").Value = "Data Copy "

$data = @(
    @(327736706,   1909893371, 1908259201, 1781589540),
    @(3007405114,  2582338660, 1752484104, 1981037516),
    @(2001113863,  2078005921, 1599648095, 1675578076),
    @(5738582713,  1991966214, 1584572805, 2259836734),
    @(1440270499,  1723703388, 1612499811, 1894197717),
    @(1281181753,  1638359940, 2904203229, 2011930042),
    @(2111512815,  2371708251, 2116357931, 4389688270),
    @(2074325641,  2297538173, 2006908595, 4113563457),
    @(2130233647,  1456606153, 1624242950, 2921578741),
    @(2117753470,  1725827280, 2629873057, 2207160820)
)
for ($i = 0; $i -lt 10; $i++) {
    $r = 3 + $i
    for ($j = 0; $j -lt 4; $j++) {
        $c = 2 + $j
        $ws.Cells.Item($r, $c).Value = $data[$i][$j]
    }
}

$ws.Range("B13").Formula = "=AVERAGE(B3:B12)"
$ws.Range("C13").Formula = "=AVERAGE(C3:C12)"
$ws.Range("D13").Formula = "=AVERAGE(D3:D12)"
$ws.Range("E13").Formula = "=AVERAGE(E3:E12)"

# Give the average row the same bold/centered look used on A13.
$ws.Range("A13").Copy()
$ws.Range("B13:E13").PasteSpecial(-4122)

[void]$ws.Range("G12").Select()

# ---------------------------------------------------------------------
# 2. Same average-row restyle on the two sheets whose formulas were
#    reformatted alongside the new data.
# ---------------------------------------------------------------------
$javaSockets = $wb.Worksheets.Item("JAVA SOCKETS RESULTS")
$javaSockets.Range("A13").Copy()
$javaSockets.Range("B13:E13").PasteSpecial(-4122)

$memMap.Range("A13").Copy()
$memMap.Range("B13:E13").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Selection / active-tab bookkeeping.
# ---------------------------------------------------------------------
$javaSockets.Activate()
[void]$javaSockets.Range("B13:E13").Select()

$memMap.Activate()
[void]$memMap.Range("A13:E13").Select()

$msgSent = $wb.Worksheets.Item("1 Message Sent")
$msgSent.Activate()
[void]$msgSent.Range("K15").Select()
